$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for specific rows per repulled data
$ws.Range("F5").Value = 5
$ws.Range("F6").Value = 2
$ws.Range("F7").Value = -2
$ws.Range("F9").Value = -6
$ws.Range("F11").Value = 1
$ws.Range("F12").Value = -8
